# Append a new data row (row 91) to the bottom of the trimmed graph data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A91").Value = 1.269
$ws.Range("B91").Value = 1.52
$ws.Range("C91").Value = 2.263
